# Applies the content edits described by the commit "Add files via upload":
#   - Slide 3 (sldId 273): tweak the "Texas / California" sales bullets.
#   - Slide 4 (sldId 288): tweak title + reorder/rebold the top-3-states bullets.
#   - Slide 5 (sldId 289): split the discounts sentence into two runs.
#   - Slide 7 (sldId 293): add a trailing period after "EasyPay".
#
# Helper: PowerPoint's TextRange.Text setter diffs the new string against the
# old one and only rewrites the runs that actually changed (so it tends to
# split a run right where the characters start to differ, which isn't always
# where we want the bold/non-bold boundary to land). Setting the paragraph to
# a disjoint placeholder string first means the final assignment has nothing
# in common with the previous text, so it always lands back in a single run
# with the base run's formatting - giving us a clean slate to re-split with
# Characters(start,length) + explicit Font formatting.

$p = $ppt.ActivePresentation

function Set-ParaTextClean {
    param($para, [string]$newText)
    $para.Text = "`u{E000}`u{E000}`u{E000}`u{E000}`u{E000}`u{E000}`u{E000}`u{E000}"
    $para.Text = $newText
}

# ---------------------------------------------------------------------------
# Slide 3 (sldId 273) - shape 18 "Text Placeholder 17"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sp3 = $s3.Shapes.Item(3)
$tr3 = $sp3.TextFrame.TextRange

# Paragraph 3: "California has second highest sales in US"
#           -> "California has second highest sales " (bold) + "at $15.39 Million" (not bold)
$para3_3 = $tr3.Paragraphs(3, 1)
Set-ParaTextClean $para3_3 "California has second highest sales at `$15.39 Million"
$boldPrefix = "California has second highest sales "
$suffix = $para3_3.Text.Substring($boldPrefix.Length)
$suffixRange = $para3_3.Characters($boldPrefix.Length + 1, $suffix.Length)
$suffixRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# Slide 4 (sldId 288) - shape 17 "Title 16" and shape 18 "Text Placeholder 17"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$title4 = $s4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "What are the top 3 US States with highest orders per month?"

$body4 = $s4.Shapes.Item(2)
$tr4 = $body4.TextFrame.TextRange

# Paragraph 1: "Texas - 1 (50.91K Sales)" -> "1 - Texas " (bold) + "(50.91K Sales)" (not bold)
$para4_1 = $tr4.Paragraphs(1, 1)
Set-ParaTextClean $para4_1 "1 - Texas (50.91K Sales)"
$prefix = "1 - Texas "
$rest = $para4_1.Text.Substring($prefix.Length)
$restRange = $para4_1.Characters($prefix.Length + 1, $rest.Length)
$restRange.Font.Bold = $false

# Paragraph 2: "California - 2 (49.83K Sales)" -> "2 - California " (bold) + "(49.83K Sales)" (not bold)
$para4_2 = $tr4.Paragraphs(2, 1)
Set-ParaTextClean $para4_2 "2 - California (49.83K Sales)"
$prefix = "2 - California "
$rest = $para4_2.Text.Substring($prefix.Length)
$restRange = $para4_2.Characters($prefix.Length + 1, $rest.Length)
$restRange.Font.Bold = $false

# Paragraph 3: "New York - 3 (48.07K Sales)" -> "3 - New York " (bold) + "(48.07K Sales)" (not bold)
$para4_3 = $tr4.Paragraphs(3, 1)
Set-ParaTextClean $para4_3 "3 - New York (48.07K Sales)"
$prefix = "3 - New York "
$rest = $para4_3.Text.Substring($prefix.Length)
$restRange = $para4_3.Characters($prefix.Length + 1, $rest.Length)
$restRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# Slide 5 (sldId 289) - shape 18 "Text Placeholder 17"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$body5 = $s5.Shapes.Item(3)
$tr5 = $body5.TextFrame.TextRange

# Single bold run -> "Based on ... from " (not bold) + "November ... discounts." (bold)
$para5_1 = $tr5.Paragraphs(1, 1)
$notBoldPrefix = "Based on the visualization, the discounts start to increase from "
$prefixRange = $para5_1.Characters(1, $notBoldPrefix.Length)
$prefixRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# Slide 7 (sldId 293) - shape 18 "Text Placeholder 17"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$body7 = $s7.Shapes.Item(3)
$tr7 = $body7.TextFrame.TextRange

# Paragraph 3: "It is followed by EasyPay" -> "It is followed by EasyPay."
# (appended via InsertAfter so the existing "It is followed by "/"EasyPay"
# run split - and EasyPay's err="1" flag - survive; only a new run is added
# for the trailing period.)
$para7_3 = $tr7.Paragraphs(3, 1)
$null = $para7_3.InsertAfter(".")
$dotRange = $para7_3.Characters($para7_3.Text.Length, 1)
$dotRange.Font.Bold = $true
